$wb = $excel.ActiveWorkbook

# Sheet 1: TC02_Verify_MYACC_RegisteredUse
$ws1 = $wb.Worksheets.Item("TC02_Verify_MYACC_RegisteredUse")
# Sheet 2: Testdata
$ws2 = $wb.Worksheets.Item("Testdata")

# Order of writes matters for shared-string table ordering: write in the
# same first-appearance order as the target workbook (53..56):
#   53 PurchasingHistory, 54 Purchasing History, 55 QuickOrderMyacc, 56 Storerooms
$ws1.Range("C18").Value = "PurchasingHistory"
$ws2.Range("B14").Value = "Purchasing History"
$ws1.Range("C19").Value = "QuickOrderMyacc"
$ws1.Range("C20").Value = "Storerooms"
$ws2.Range("B16").Value = "Storerooms"

# View state changes: re-create the final selection/active-sheet state.
# Sheet1 (TC02...) is no longer the active tab, and is scrolled so row 7 is
# at the top with C25 selected.
$ws1.Activate()
$ws1.Range("C25").Select()
$excel.ActiveWindow.ScrollRow = 7

# Sheet2 (Testdata) ends up active/selected, scrolled so row 12 is at the
# top with B19 selected.
$ws2.Activate()
$ws2.Range("B19").Select()
$excel.ActiveWindow.ScrollRow = 12
